# C5-PowerPoint.pptx edit
#
# 1) The table on slide 6 gets a new (built-in) table style applied.
# 2) The deck's theme palette is reset from the custom "Integral" colour
#    set back to the standard "Office Theme" colour set.

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{221CB12D-C1AE-44EB-9294-1A0CCAB75340}")
    }
}

# --- 2) Theme colours -------------------------------------------------
# Restore the stock "Office Theme" colour scheme (currently the deck
# carries the custom "Integral" palette).
$firstSlide = $p.Slides.Item(1)
$colors = $firstSlide.ThemeColorScheme

# PowerPoint's RGB() long is 0x00BBGGRR, i.e. R + G*256 + B*65536.
$colors.Item(1).RGB  = 0x000000   # Dark 1      #000000
$colors.Item(2).RGB  = 0xFFFFFF   # Light 1     #FFFFFF
$colors.Item(3).RGB  = 0x6A5444   # Dark 2      #44546A
$colors.Item(4).RGB  = 0xE6E6E7   # Light 2     #E7E6E6
$colors.Item(5).RGB  = 0xD59B5B   # Accent 1    #5B9BD5
$colors.Item(6).RGB  = 0x317DED   # Accent 2    #ED7D31
$colors.Item(7).RGB  = 0xA5A5A5   # Accent 3    #A5A5A5
$colors.Item(8).RGB  = 0x00C0FF   # Accent 4    #FFC000
$colors.Item(9).RGB  = 0xC47244   # Accent 5    #4472C4
$colors.Item(10).RGB = 0x47AD70   # Accent 6    #70AD47
$colors.Item(11).RGB = 0xC16305   # Hyperlink   #0563C1
$colors.Item(12).RGB = 0x724F95   # Followed hyperlink #954F72
